$wb = $excel.ActiveWorkbook

# --- Update template version on the "isa_template" sheet ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.3"

# --- Update building block headers / term accession URLs on the "New Table" sheet ---
$wsTable = $wb.Worksheets.Item("New Table")

# Rename building blocks from Parameter to Component
$wsTable.Range("H1").Value = "Component [NMR sample tube]"
$wsTable.Range("K1").Value = "Component [NMR solvent]"

# Update term accession number URLs to the new bioregistry.io format
$wsTable.Range("D2").Value = "https://bioregistry.io/OBI:0000516"
$wsTable.Range("G2").Value = "https://bioregistry.io/CHEBI:17790"
$wsTable.Range("M2").Value = "https://bioregistry.io/NCIT:C91099"
$wsTable.Range("T2").Value = "https://bioregistry.io/UO:0000027"

# Keep the underlying Excel table's column names in sync with the header cells
$table = $wsTable.ListObjects.Item("annotationTable")
$table.ListColumns.Item(8).Name = "Component [NMR sample tube]"
$table.ListColumns.Item(11).Name = "Component [NMR solvent]"
